$d = $word.ActiveDocument

$d.Content.Find.Execute("396÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "869÷6=", 2) | Out-Null
$d.Content.Find.Execute("274÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "236÷6=", 2) | Out-Null
$d.Content.Find.Execute("360÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "355÷8=", 2) | Out-Null
$d.Content.Find.Execute("921÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "994÷2=", 2) | Out-Null
$d.Content.Find.Execute("315÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "533÷5=", 2) | Out-Null
$d.Content.Find.Execute("734÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "328÷7=", 2) | Out-Null
$d.Content.Find.Execute("582÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "195÷5=", 2) | Out-Null
$d.Content.Find.Execute("515÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "948÷9=", 2) | Out-Null
$d.Content.Find.Execute("914÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "920÷4=", 2) | Out-Null
$d.Content.Find.Execute("270÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "458÷2=", 2) | Out-Null
$d.Content.Find.Execute("557÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "452÷2=", 2) | Out-Null
$d.Content.Find.Execute("755÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "900÷2=", 2) | Out-Null
$d.Content.Find.Execute("468÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "429÷8=", 2) | Out-Null
$d.Content.Find.Execute("220÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "107÷3=", 2) | Out-Null
$d.Content.Find.Execute("573÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "216÷7=", 2) | Out-Null
$d.Content.Find.Execute("993÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "190÷5=", 2) | Out-Null
$d.Content.Find.Execute("930÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "410÷6=", 2) | Out-Null
$d.Content.Find.Execute("680÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "721÷7=", 2) | Out-Null
$d.Content.Find.Execute("930÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "762÷9=", 2) | Out-Null
$d.Content.Find.Execute("478÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "238÷5=", 2) | Out-Null
$d.Content.Find.Execute("790÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "679÷8=", 2) | Out-Null
$d.Content.Find.Execute("491÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "287÷5=", 2) | Out-Null
$d.Content.Find.Execute("723÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "752÷9=", 2) | Out-Null
$d.Content.Find.Execute("380÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "239÷8=", 2) | Out-Null
$d.Content.Find.Execute("759÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "443÷8=", 2) | Out-Null
